$d = $word.ActiveDocument

$d.Content.Find.Execute("49×14=686", $true, $false, $false, $false, $false, $true, 1, $false, "95×84=7980", 2) | Out-Null
$d.Content.Find.Execute("49×36=1764", $true, $false, $false, $false, $false, $true, 1, $false, "97×92=8924", 2) | Out-Null
$d.Content.Find.Execute("29×12=348", $true, $false, $false, $false, $false, $true, 1, $false, "98×58=5684", 2) | Out-Null
$d.Content.Find.Execute("89×77=6853", $true, $false, $false, $false, $false, $true, 1, $false, "61×79=4819", 2) | Out-Null
$d.Content.Find.Execute("75×42=3150", $true, $false, $false, $false, $false, $true, 1, $false, "30×85=2550", 2) | Out-Null
$d.Content.Find.Execute("59×53=3127", $true, $false, $false, $false, $false, $true, 1, $false, "25×29=725", 2) | Out-Null
$d.Content.Find.Execute("54×92=4968", $true, $false, $false, $false, $false, $true, 1, $false, "20×99=1980", 2) | Out-Null
$d.Content.Find.Execute("16×38=608", $true, $false, $false, $false, $false, $true, 1, $false, "67×45=3015", 2) | Out-Null
$d.Content.Find.Execute("33×94=3102", $true, $false, $false, $false, $false, $true, 1, $false, "32×28=896", 2) | Out-Null
$d.Content.Find.Execute("29×46=1334", $true, $false, $false, $false, $false, $true, 1, $false, "24×84=2016", 2) | Out-Null
$d.Content.Find.Execute("15×27=405", $true, $false, $false, $false, $false, $true, 1, $false, "39×90=3510", 2) | Out-Null
$d.Content.Find.Execute("12×13=156", $true, $false, $false, $false, $false, $true, 1, $false, "41×36=1476", 2) | Out-Null
$d.Content.Find.Execute("36×39=1404", $true, $false, $false, $false, $false, $true, 1, $false, "86×86=7396", 2) | Out-Null
$d.Content.Find.Execute("25×68=1700", $true, $false, $false, $false, $false, $true, 1, $false, "74×85=6290", 2) | Out-Null
$d.Content.Find.Execute("27×73=1971", $true, $false, $false, $false, $false, $true, 1, $false, "86×68=5848", 2) | Out-Null
$d.Content.Find.Execute("45×91=4095", $true, $false, $false, $false, $false, $true, 1, $false, "63×47=2961", 2) | Out-Null
$d.Content.Find.Execute("36×85=3060", $true, $false, $false, $false, $false, $true, 1, $false, "13×68=884", 2) | Out-Null
$d.Content.Find.Execute("67×25=1675", $true, $false, $false, $false, $false, $true, 1, $false, "91×46=4186", 2) | Out-Null
$d.Content.Find.Execute("79×71=5609", $true, $false, $false, $false, $false, $true, 1, $false, "89×50=4450", 2) | Out-Null
$d.Content.Find.Execute("32×45=1440", $true, $false, $false, $false, $false, $true, 1, $false, "61×76=4636", 2) | Out-Null
$d.Content.Find.Execute("52×73=3796", $true, $false, $false, $false, $false, $true, 1, $false, "74×73=5402", 2) | Out-Null
$d.Content.Find.Execute("96×78=7488", $true, $false, $false, $false, $false, $true, 1, $false, "57×48=2736", 2) | Out-Null
$d.Content.Find.Execute("15×53=795", $true, $false, $false, $false, $false, $true, 1, $false, "70×22=1540", 2) | Out-Null
$d.Content.Find.Execute("64×30=1920", $true, $false, $false, $false, $false, $true, 1, $false, "73×34=2482", 2) | Out-Null
$d.Content.Find.Execute("52×61=3172", $true, $false, $false, $false, $false, $true, 1, $false, "27×29=783", 2) | Out-Null
